$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-08-18 08:16:51"
$overview.Range("G5").Value = "2016-08-18 08:16:51"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("E5").Value = "mt"
$zhcn.Range("H4").Value = "2016-08-18 08:16:45"
$zhcn.Range("H5").Value = "2016-08-18 08:16:45"
$zhcn.Range("K4").Value = "2016-08-18 08:17:19"
$zhcn.Range("K5").Value = "2016-08-18 08:17:19"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "mt"
$dede.Range("E5").Value = "mt"
$dede.Range("K4").Value = "2016-08-18 08:17:27"
$dede.Range("K5").Value = "2016-08-18 08:17:27"
